$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.930.48'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.760.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.95'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3761'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3345'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.67'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.117'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07141'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -6.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.24'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.170'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -5.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.135'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.756.95'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001049'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06568'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '79.99'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.34%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.247'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.930.52'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.65'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -8.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.382'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.09'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.71'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -8.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.314'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -10.45%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.266'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -15.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.38'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.022'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.762'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -7.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08733'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.13'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -8.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02327'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.03%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -6.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06165'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.124'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -7.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2101'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.89%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.451'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -10.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.004'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.71'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.31%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5997'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.93'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.999'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -8.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07157'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.177'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.66%  '
